# "add ID species from plate PocHistone 015"
#
# PocHistone RLFP 015 plate is now fully scored, so:
#   1. mark that sheet DONE (rename it)
#   2. the species-ID row (row 11, "Plate needed" + per-well scores) can now
#      be filled in on the two downstream plates (017 and 018) that were
#      waiting on those samples.

$wb = $excel.ActiveWorkbook

# 1) PocHistone RLFP 015 is complete -> prefix with "DONE "
$ws015 = $wb.Worksheets.Item("PocHistone RLFP 015")
$ws015.Name = "DONE PocHistone RLFP 015"

# 2) PocHistone RLFP 017 -> add row 11 (ID species / plate-needed scores)
$ws017 = $wb.Worksheets.Item("PocHistone RLFP 017")
$ws017.Range("A11").Value = "Plate needed"
$ws017.Range("B11").Value = 20
$ws017.Range("C11").Value = 31
$ws017.Range("D11").Value = 23
$ws017.Range("E11").Value = 30
$ws017.Range("F11").Value = 27
$ws017.Range("A1:M11").Select() | Out-Null

# 3) PocHistone RLFP 018 -> add row 11 (ID species / plate-needed scores)
$ws018 = $wb.Worksheets.Item("PocHistone RLFP 018")
$ws018.Range("A11").Value = "Plate needed"
$ws018.Range("B11").Value = 28
$ws018.Range("C11").Value = 27
$ws018.Range("D11").Value = 23
$ws018.Range("E11").Value = 31
$ws018.Range("F11").Value = 28
$ws018.Range("G11").Value = 24
$ws018.Range("H11").Value = 12
$ws018.Range("I11").Value = 15
$ws018.Range("J11").Value = 38
$ws018.Range("A1:M11").Select() | Out-Null

# leave the workbook focused back on the plate that was just completed
$ws015.Activate()
